# Week 15 simulations - update Rushing and Receiving stats for Broncos players

$wb = $excel.ActiveWorkbook

# ----- Rushing sheet -----
$rushing = $wb.Worksheets.Item("Rushing")

# Row 3 - M.Gordon
$rushing.Range("C3").Value = 83
$rushing.Range("D3").Value = 62
$rushing.Range("E3").Value = 10
$rushing.Range("F3").Value = 33

# Row 4 - J.Williams
$rushing.Range("C4").Value = 82
$rushing.Range("D4").Value = 55
$rushing.Range("E4").Value = 17
$rushing.Range("F4").Value = 21

# ----- Receiving sheet -----
$receiving = $wb.Worksheets.Item("Receiving")

# Row 3 - J.Williams
$receiving.Range("C3").Value = 42
$receiving.Range("D3").Value = 32
$receiving.Range("G3").Value = 7
$receiving.Range("H3").Value = 5

# Row 5 - C.Sutton
$receiving.Range("C5").Value = 52
$receiving.Range("D5").Value = 40
$receiving.Range("E5").Value = 26

# Row 6 - J.Jeudy
$receiving.Range("C6").Value = 70
$receiving.Range("D6").Value = 57
$receiving.Range("E6").Value = 20
$receiving.Range("F6").Value = 16

# Row 7 - T.Patrick
$receiving.Range("C7").Value = 52
$receiving.Range("D7").Value = 35
$receiving.Range("E7").Value = 14
$receiving.Range("G7").Value = 10
$receiving.Range("H7").Value = 4

# Row 8 - D.Spencer
$receiving.Range("C8").Value = 4

# Row 10 - N.Fant
$receiving.Range("C10").Value = 66
$receiving.Range("D10").Value = 53
$receiving.Range("E10").Value = 11
$receiving.Range("F10").Value = 6

# Row 11 - A.Okwuegbunam
$receiving.Range("C11").Value = 30
$receiving.Range("D11").Value = 28
$receiving.Range("G11").Value = 4
$receiving.Range("H11").Value = 3
